# Auto-generated edit script: refresh Ravana_Profits Sheets market-price derived values
# (currentAveragePrice / LevePrice / LeveProfit columns) per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 1104.8235
$ws.Range("I15").Value2 = 1104.8235
$ws.Range("K15").Value2 = 3314.4705
$ws.Range("M15").Value2 = -3145.4705
$ws.Range("H43").Value2 = 7024.75
$ws.Range("I43").Value2 = 6349.5
$ws.Range("K43").Value2 = 6349.5
$ws.Range("M43").Value2 = -6280.5
$ws.Range("H92").Value2 = 1192
$ws.Range("I92").Value2 = 1192
$ws.Range("K92").Value2 = 1192
$ws.Range("M92").Value2 = 56
$ws.Range("H132").Value2 = 969.5
$ws.Range("I132").Value2 = 969.5
$ws.Range("K132").Value2 = 2908.5
$ws.Range("M132").Value2 = -378.5
$ws.Range("H137").Value2 = 3636.0688
$ws.Range("I137").Value2 = 658.7273
$ws.Range("K137").Value2 = 1976.1819
$ws.Range("M137").Value2 = 573.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5331.64
$ws.Range("I32").Value2 = 4273.5654
$ws.Range("K32").Value2 = 4273.5654
$ws.Range("M32").Value2 = -3986.5654
$ws.Range("H61").Value2 = 2753.8
$ws.Range("I61").Value2 = 2542.5
$ws.Range("J61").Value2 = 3599
$ws.Range("K61").Value2 = 2542.5
$ws.Range("L61").Value2 = 3599
$ws.Range("M61").Value2 = -2330.5
$ws.Range("N61").Value2 = -4023
$ws.Range("H74").Value2 = 8498.083000000001
$ws.Range("I74").Value2 = 5747.375
$ws.Range("J74").Value2 = 13999.5
$ws.Range("K74").Value2 = 5747.375
$ws.Range("L74").Value2 = 13999.5
$ws.Range("M74").Value2 = -4873.375
$ws.Range("N74").Value2 = -15747.5
$ws.Range("H77").Value2 = 8498.083000000001
$ws.Range("I77").Value2 = 5747.375
$ws.Range("J77").Value2 = 13999.5
$ws.Range("K77").Value2 = 28736.875
$ws.Range("L77").Value2 = 69997.5
$ws.Range("M77").Value2 = -24368.875
$ws.Range("N77").Value2 = -78733.5
$ws.Range("H97").Value2 = 499.41666
$ws.Range("I97").Value2 = 499.41666
$ws.Range("K97").Value2 = 499.41666
$ws.Range("M97").Value2 = -3.416659999999979
$ws.Range("H132").Value2 = 3340.524
$ws.Range("I132").Value2 = 2550.923
$ws.Range("J132").Value2 = 4623.625
$ws.Range("K132").Value2 = 7652.768999999999
$ws.Range("L132").Value2 = 13870.875
$ws.Range("M132").Value2 = -5122.768999999999
$ws.Range("N132").Value2 = -18930.875
$ws.Range("H136").Value2 = 2753.8
$ws.Range("I136").Value2 = 2542.5
$ws.Range("J136").Value2 = 3599
$ws.Range("K136").Value2 = 7627.5
$ws.Range("L136").Value2 = 10797
$ws.Range("M136").Value2 = -5077.5
$ws.Range("N136").Value2 = -15897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value2 = 2535.5
$ws.Range("I105").Value2 = 2428.8333
$ws.Range("K105").Value2 = 2428.8333
$ws.Range("M105").Value2 = -681.8332999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 90476640
$ws.Range("J4").Value2 = 133511440
$ws.Range("L4").Value2 = 400534320
$ws.Range("N4").Value2 = -400534544
$ws.Range("H98").Value2 = 3629.6667
$ws.Range("J98").Value2 = 4177.6
$ws.Range("L98").Value2 = 12532.8
$ws.Range("N98").Value2 = -15528.8
$ws.Range("H122").Value2 = 3724.4443
$ws.Range("J122").Value2 = 3816.3022
$ws.Range("L122").Value2 = 34346.7198
$ws.Range("N122").Value2 = -39246.7198
$ws.Range("H132").Value2 = 3858.75
$ws.Range("I132").Value2 = 966
$ws.Range("J132").Value2 = 5594.4
$ws.Range("K132").Value2 = 8694
$ws.Range("L132").Value2 = 50349.6
$ws.Range("M132").Value2 = -6164
$ws.Range("N132").Value2 = -55409.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 3373.7144
$ws.Range("I80").Value2 = 2334.818
$ws.Range("J80").Value2 = 4516.5
$ws.Range("K80").Value2 = 2334.818
$ws.Range("L80").Value2 = 4516.5
$ws.Range("M80").Value2 = -1336.818
$ws.Range("N80").Value2 = -6512.5
$ws.Range("H83").Value2 = 3373.7144
$ws.Range("I83").Value2 = 2334.818
$ws.Range("J83").Value2 = 4516.5
$ws.Range("K83").Value2 = 11674.09
$ws.Range("L83").Value2 = 22582.5
$ws.Range("M83").Value2 = -6682.09
$ws.Range("N83").Value2 = -32566.5
$ws.Range("H97").Value2 = 720.4
$ws.Range("I97").Value2 = 578.2222
$ws.Range("J97").Value2 = 2000
$ws.Range("K97").Value2 = 578.2222
$ws.Range("L97").Value2 = 2000
$ws.Range("M97").Value2 = -82.22220000000004
$ws.Range("N97").Value2 = -2992
$ws.Range("H102").Value2 = 3458.8
$ws.Range("I102").Value2 = 3458.8
$ws.Range("K102").Value2 = 3458.8
$ws.Range("M102").Value2 = -1836.8
$ws.Range("H122").Value2 = 1140.5714
$ws.Range("I122").Value2 = 997.4
$ws.Range("J122").Value2 = 1498.5
$ws.Range("K122").Value2 = 2992.2
$ws.Range("L122").Value2 = 4495.5
$ws.Range("M122").Value2 = -542.1999999999998
$ws.Range("N122").Value2 = -9395.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4249.5
$ws.Range("J7").Value2 = 4499.5
$ws.Range("L7").Value2 = 4499.5
$ws.Range("N7").Value2 = -4723.5
$ws.Range("H16").Value2 = 1054.2
$ws.Range("I16").Value2 = 1193.5
$ws.Range("K16").Value2 = 1193.5
$ws.Range("M16").Value2 = -1023.5
$ws.Range("H22").Value2 = 3932.8215
$ws.Range("I22").Value2 = 3894.6
$ws.Range("K22").Value2 = 3894.6
$ws.Range("M22").Value2 = -3599.6
$ws.Range("H27").Value2 = 3932.8215
$ws.Range("I27").Value2 = 3894.6
$ws.Range("K27").Value2 = 3894.6
$ws.Range("M27").Value2 = -3787.6
$ws.Range("H40").Value2 = 4966.3335
$ws.Range("I40").Value2 = 4997
$ws.Range("K40").Value2 = 4997
$ws.Range("M40").Value2 = -4861
$ws.Range("H68").Value2 = 99999
$ws.Range("I68").Value2 = 99999
$ws.Range("K68").Value2 = 99999
$ws.Range("M68").Value2 = -99250
$ws.Range("H71").Value2 = 99999
$ws.Range("I71").Value2 = 99999
$ws.Range("K71").Value2 = 499995
$ws.Range("M71").Value2 = -496251
$ws.Range("H126").Value2 = 4249.5
$ws.Range("J126").Value2 = 4499.5
$ws.Range("L126").Value2 = 13498.5
$ws.Range("N126").Value2 = -18438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 2996.5
$ws.Range("I122").Value2 = 2662.3333
$ws.Range("K122").Value2 = 7986.999899999999
$ws.Range("M122").Value2 = -5536.999899999999
$ws.Range("H132").Value2 = 2886.2727
$ws.Range("I132").Value2 = 2583.6296
$ws.Range("J132").Value2 = 4248.1665
$ws.Range("K132").Value2 = 7750.888800000001
$ws.Range("L132").Value2 = 12744.4995
$ws.Range("M132").Value2 = -5220.888800000001
$ws.Range("N132").Value2 = -17804.4995

Write-Output "Updated 161 cells across 7 sheets"